$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: was empty, now gets the "Dynamic API requests..." note (written first so it
# lands at shared-string index 33, matching the target ordering)
$ws.Range("E27").Value = "Dynamic API requests and connecting input fields to character fetch. Dealing with problems of passing data from stless/stful widget"

# New "Hours:" column F labels on the weekly summary rows (written next, so "Hours:"
# becomes shared-string index 34)
$ws.Range("F4").Value = "Hours:"
$ws.Range("F9").Value = "Hours:"
$ws.Range("F14").Value = "Hours:"
$ws.Range("F19").Value = "Hours:"
$ws.Range("F24").Value = "Hours:"
$ws.Range("F29").Value = "Hours:"
$ws.Range("F34").Value = "Hours:"

# Row 24 description updated to mention the form as well (written last, landing at
# shared-string index 35)
$ws.Range("E24").Value = "Created EU realm list, search page & form"

# Update the view: scroll back to the top-left and move the active selection to J21
$ws.Range("J21").Select()
